# Update "想去人数" (F column) values on the sheets that contain this
# event table: "展览" (sheet1) and "全部类型" (sheet4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 6752
    3 = 47
    4 = 195
    5 = 1056
    6 = 146
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
